$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date formatting (style) from the last existing date cell (A101)
# down onto the new date cells A102:A110 so they reuse the workbook's
# built-in date style instead of minting a new custom number format.
$ws.Range("A101").Copy() | Out-Null
$ws.Range("A102:A110").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$rows = @(
    @{ Row = 102; Date = 43224; B = "2018-05-04 1.JPG";  C = "Die ESCs sind angekommen. Der Neue ist etwas größer, kam aber schneller und ist auch für eventuelle stärkere Motoren geeignet" },
    @{ Row = 103; Date = 43224; B = "2018-05-04 2.JPG";  C = "Die Teststation ist aufgebaut. Die ersten Flugversuche können gestartet werden" },
    @{ Row = 104; Date = 43224; B = "2018-05-04 3.JPG";  C = "Die Teststation von Oben" },
    @{ Row = 105; Date = 43224; B = "2018-05-04 4.AVI";  C = "Das erste mal drehen sich die Motoren mit en neuen Platinen, ESCs und Programmen. Auf Anhieb wohl gemerkt !!" },
    @{ Row = 106; Date = 43224; B = "2018-05-04 5.AVI";  C = "Der erste Flug aus Sicht der Steuerung. Ausgabe der Drohne, des Arduinos und das GUI" },
    @{ Row = 107; Date = 43224; B = "2018-05-04 6.AVI";  C = "Zögerliches Anheben der Rotoren" },
    @{ Row = 108; Date = 43224; B = "2018-05-04 7.AVI";  C = "CRASH! Etliche Fehlversuche später und lösen einiger Verbindungsprobleme resultierte dan in diesem Test. Erwartungsgemäß verhielt sich die nicht-kalibirerte Drohne seltsam und zerstörte sich selbst. Wenigstens sind die Motoren jetzt sicher stark genug" },
    @{ Row = 109; Date = 43224; B = "2018-05-04 8.mp4";  C = "Der Crash aus sicht des Basiscomputers mit anfänglichen Startschwierigkeiten und meinen Erklärenden Selbstgesprächen (was es sein soll, weiß ich selbst nicht)" },
    @{ Row = 110; Date = 43224; B = "2018-05-04 9.avi";  C = "Wir sind jetzt auf YouTube. Der erste Clip wurde zur Demonstration zusammengeschnitten und hochgeladen." }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value2 = $r.Date
    $ws.Cells.Item($r.Row, 2).Value2 = $r.B
    $ws.Cells.Item($r.Row, 3).Value2 = $r.C
}

# Update the visible selection to match the final state (last edited cell).
$ws.Range("C110").Select() | Out-Null
